# eprepago.xlsx — add a "correoUsuario" column (N) with a mailto hyperlink
# for the new e-prepago enrollment test row (per commit: "Se agrega correo
# en el archivo datadriven, inscripcion de e-prepago").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# New header in N1.
$ws.Range("N1").Value = "correoUsuario"

# New value in N2 — set the display text first so the hyperlink insert
# below reuses it verbatim instead of writing the raw mailto: address.
$ws.Range("N2").Value = "jruav@devco.com.co"

# Turn N2 into a mailto hyperlink (Insert > Hyperlink > E-mail Address).
$null = $ws.Hyperlinks.Add($ws.Range("N2"), "mailto:jruav@devco.com.co")

# Match the selection left behind by inserting the hyperlink on N2.
$null = $ws.Range("N2").Select()
